# Auto-generated edit script: apply scheduled-runner market-price/profit updates
# to the Kraken_Profits workbook (FFXIV leve-profit tracker).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 53
$ws.Range("H53").Value = 1088.375
$ws.Range("I53").Value = 101.166664
$ws.Range("K53").Value = 101.166664
$ws.Range("M53").Value = 535.833336
# Row 62
$ws.Range("H62").Value = 9583.333000000001
$ws.Range("I62").Value = 9675
$ws.Range("K62").Value = 9675
$ws.Range("M62").Value = -9051
# Row 65
$ws.Range("H65").Value = 9583.333000000001
$ws.Range("I65").Value = 9675
$ws.Range("K65").Value = 48375
$ws.Range("M65").Value = -45255
# Row 113
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -7508

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3087.2144
$ws.Range("I32").Value = 2690.077
$ws.Range("J32").Value = 8250
$ws.Range("K32").Value = 2690.077
$ws.Range("L32").Value = 8250
$ws.Range("M32").Value = -2403.077
$ws.Range("N32").Value = -8824
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 45
$ws.Range("H45").Value = 2719.4783
$ws.Range("I45").Value = 2459.6667
$ws.Range("K45").Value = 2459.6667
$ws.Range("M45").Value = -2082.6667
# Row 76
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10676
# Row 79
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12340
# Row 97
$ws.Range("H97").Value = 1517.9231
$ws.Range("I97").Value = 703.7778
$ws.Range("J97").Value = 3349.75
$ws.Range("K97").Value = 703.7778
$ws.Range("L97").Value = 3349.75
$ws.Range("M97").Value = -207.7778
$ws.Range("N97").Value = -4341.75
# Row 110
$ws.Range("H110").Value = 3954.7
$ws.Range("I110").Value = 574.5
$ws.Range("J110").Value = 4799.75
$ws.Range("K110").Value = 574.5
$ws.Range("L110").Value = 4799.75
$ws.Range("M110").Value = 1470.5
$ws.Range("N110").Value = -8889.75
# Row 132
$ws.Range("H132").Value = 2424.5715
$ws.Range("I132").Value = 2009.6666
$ws.Range("K132").Value = 6028.9998
$ws.Range("M132").Value = -3498.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3548.25
$ws.Range("I86").Value = 1937.2
$ws.Range("J86").Value = 6233.3335
$ws.Range("K86").Value = 1937.2
$ws.Range("L86").Value = 6233.3335
$ws.Range("M86").Value = -814.2
$ws.Range("N86").Value = -8479.333500000001
# Row 89
$ws.Range("H89").Value = 3548.25
$ws.Range("I89").Value = 1937.2
$ws.Range("J89").Value = 6233.3335
$ws.Range("K89").Value = 9686
$ws.Range("L89").Value = 31166.6675
$ws.Range("M89").Value = -4070
$ws.Range("N89").Value = -42398.6675
# Row 107
$ws.Range("H107").Value = 22166.5
$ws.Range("I107").Value = 1230.5
$ws.Range("J107").Value = 43102.5
$ws.Range("K107").Value = 1230.5
$ws.Range("L107").Value = 43102.5
$ws.Range("M107").Value = 689.5
$ws.Range("N107").Value = -46942.5
# Row 134
$ws.Range("H134").Value = 4949.75
$ws.Range("I134").Value = 2599.6667
$ws.Range("K134").Value = 7799.000100000001
$ws.Range("M134").Value = -5264.000100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 4436
$ws.Range("I86").Value = 4436
$ws.Range("K86").Value = 4436
$ws.Range("M86").Value = -3313
# Row 89
$ws.Range("H89").Value = 4436
$ws.Range("I89").Value = 4436
$ws.Range("K89").Value = 22180
$ws.Range("M89").Value = -16564
# Row 102
$ws.Range("H102").Value = 45000
$ws.Range("J102").Value = 45000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -49868
# Row 107
$ws.Range("H107").Value = 192.08333
$ws.Range("I107").Value = 111.7
$ws.Range("J107").Value = 594
$ws.Range("K107").Value = 111.7
$ws.Range("L107").Value = 594
$ws.Range("M107").Value = 1808.3
$ws.Range("N107").Value = -4434
# Row 132
$ws.Range("H132").Value = 7253.6665
$ws.Range("I132").Value = 7144.4
$ws.Range("J132").Value = 7800
$ws.Range("K132").Value = 21433.2
$ws.Range("L132").Value = 23400
$ws.Range("M132").Value = -18903.2
$ws.Range("N132").Value = -28460
# Row 134
$ws.Range("H134").Value = 2375
$ws.Range("I134").Value = 2375
$ws.Range("K134").Value = 7125
$ws.Range("M134").Value = -4590

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4000
# Row 34
$ws.Range("H34").Value = 1856.4286
$ws.Range("J34").Value = 2149.1667
$ws.Range("L34").Value = 6447.500100000001
$ws.Range("N34").Value = -6615.500100000001
# Row 64
$ws.Range("H64").Value = 898
$ws.Range("I64").Value = 898
$ws.Range("K64").Value = 2694
$ws.Range("M64").Value = -2424
# Row 67
$ws.Range("H67").Value = 898
$ws.Range("I67").Value = 898
$ws.Range("K67").Value = 2694
$ws.Range("M67").Value = -1758
# Row 137
$ws.Range("H137").Value = 816.6667
$ws.Range("I137").Value = 816.6667
$ws.Range("K137").Value = 2450.0001
$ws.Range("M137").Value = 2649.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 15344
$ws.Range("J95").Value = 15344
$ws.Range("L95").Value = 15344
$ws.Range("N95").Value = -20836

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5298.8
$ws.Range("I7").Value = 4831.6665
$ws.Range("K7").Value = 4831.6665
$ws.Range("M7").Value = -4719.6665
# Row 40
$ws.Range("H40").Value = 7999.6665
$ws.Range("J40").Value = 7999.5
$ws.Range("L40").Value = 7999.5
$ws.Range("N40").Value = -8271.5
# Row 53
$ws.Range("H53").Value = 46000
$ws.Range("I53").Value = 46000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 46000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -45482
$ws.Range("N53").ClearContents()
# Row 93
$ws.Range("H93").Value = 1443.9166
$ws.Range("I93").Value = 1443.9166
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1443.9166
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -195.9166
$ws.Range("N93").ClearContents()
# Row 126
$ws.Range("H126").Value = 5298.8
$ws.Range("I126").Value = 4831.6665
$ws.Range("K126").Value = 14494.9995
$ws.Range("M126").Value = -12024.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5780.8
$ws.Range("I132").Value = 5780.8
$ws.Range("K132").Value = 17342.4
$ws.Range("M132").Value = -14812.4
